$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. thousands-dot formatted prices).
# Force text format before assignment so Excel does not coerce these into numbers,
# then restore the default style so no stray formatting is introduced.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.118.21"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "1.906.31"
$ws.Range("E3").Value = "  +2.00%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").Value = "319.46"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").Value = "0.9991"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("D7").Value = "0.5050"
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").Value = "0.4071"
$ws.Range("E8").Value = "  +3.81%  "
$ws.Range("D9").Value = "0.08337"
$ws.Range("E9").Value = "  +1.90%  "
$ws.Range("D10").Value = "42.26"
$ws.Range("E10").Value = "  +0.52%  "
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("D12").Value = "23.89"
$ws.Range("E12").Value = "  +5.83%  "
$ws.Range("D13").Value = "6.388"
$ws.Range("E13").Value = "  +2.16%  "
$ws.Range("D14").Value = "1.903.61"
$ws.Range("E14").Value = "  +1.88%  "
$ws.Range("D15").Value = "7.215"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").Value = "0.9991"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").Value = "92.34"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").Value = "0.00001094"
$ws.Range("E18").Value = "  +1.60%  "
$ws.Range("D19").Value = "0.06491"
$ws.Range("E19").Value = "  +2.42%  "
$ws.Range("D20").Value = "18.26"
$ws.Range("E20").Value = "  +2.43%  "
$ws.Range("D21").Value = "0.9991"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").Value = "5.922"
$ws.Range("E22").Value = "  +2.28%  "
$ws.Range("D23").Value = "30.120.73"
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("E24").Value = "  +2.49%  "
$ws.Range("D25").Value = "2.188"
$ws.Range("E25").Value = "  -0.88%  "
$ws.Range("D26").Value = "2.120.97"
$ws.Range("D27").Value = "21.72"
$ws.Range("E27").Value = "  +4.13%  "
$ws.Range("D28").Value = "162.55"
$ws.Range("E28").Value = "  +0.77%  "
$ws.Range("E29").Value = "  +2.12%  "
$ws.Range("D31").Value = "1.139"
$ws.Range("E31").Value = "  +9.64%  "
$ws.Range("E32").Value = "  +0.79%  "
$ws.Range("D33").Value = "5.951"
$ws.Range("E33").Value = "  +1.46%  "
$ws.Range("D34").Value = "3.784"
$ws.Range("E34").Value = "  +1.35%  "
$ws.Range("D35").Value = "0.02453"
$ws.Range("E35").Value = "  +1.17%  "
$ws.Range("D36").Value = "5.354"
$ws.Range("E36").Value = "  +2.89%  "
$ws.Range("D37").Value = "0.06372"
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("D38").Value = "0.2143"
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("E39").Value = "  +4.08%  "
$ws.Range("E40").Value = "  +1.93%  "
$ws.Range("D41").Value = "8.586"
$ws.Range("E41").Value = "  +0.80%  "
$ws.Range("D42").Value = "11.37"
$ws.Range("E42").Value = "  +1.43%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "13.39"
$ws.Range("E44").Value = "  +4.98%  "
$ws.Range("D45").Value = "2.201"
$ws.Range("E45").Value = "  +10.79%  "
$ws.Range("D46").Value = "0.6073"
$ws.Range("E46").Value = "  +3.25%  "
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("D49").Value = "121.44"
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("E50").Value = "  +2.91%  "
$ws.Range("D51").Value = "1.138"
$ws.Range("E51").Value = "  +1.66%  "

$dRange.Style = "Normal"
